$wb = $excel.ActiveWorkbook

function Set-TextValue($range, $value) {
    # Write a numeric-looking string while keeping it stored as text
    # (Excel's COM layer otherwise silently coerces "0.52" etc. to a number),
    # then restore the cell to the default "Normal" style so no stray
    # number-format styling is left behind.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$totalSheet = $wb.Worksheets.Item(1)   # "总计"
$q3Sheet = $wb.Worksheets.Item(2)      # currently "2022-Q3"

# Duplicate the existing "2022-Q3" sheet right after itself, so the
# original Q3 data/formatting is preserved on its own sheet, while the
# original sheet object becomes the new "2022-Q4" sheet.
# Resulting order: 总计(1), 2022-Q4(2, was sheet2), 2022-Q3(3, new copy)
$q3Sheet.Copy($null, $q3Sheet)
$q3Sheet.Name = "2022-Q4"
$q4DataSheet = $q3Sheet

$newQ3Sheet = $wb.Worksheets.Item(3)
$newQ3Sheet.Name = "2022-Q3"

# --- Update the "总计" (totals) sheet ---
# Row 2 now reports the new 2022-Q4 totals.
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("D2").Value = 0.06

# Row 3 (new) keeps the previous 2022-Q3 totals, matching the style of row 2.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q3"
$totalSheet.Range("C3").Value = 2
$totalSheet.Range("D3").Value = 0.08

$totalSheet.Range("A2").Copy()
$totalSheet.Range("A3").PasteSpecial(-4122)

# --- Populate the "2022-Q4" sheet with the new quarter's fund data ---
Set-TextValue $q4DataSheet.Range("D2") "0.52"
Set-TextValue $q4DataSheet.Range("E2") "93.31"
Set-TextValue $q4DataSheet.Range("F2") "7.87"
Set-TextValue $q4DataSheet.Range("G2") "0.0409"
$q4DataSheet.Range("H2").Value = 6

Set-TextValue $q4DataSheet.Range("D3") "0.18"
Set-TextValue $q4DataSheet.Range("E3") "93.31"
Set-TextValue $q4DataSheet.Range("F3") "7.87"
Set-TextValue $q4DataSheet.Range("G3") "0.0142"
$q4DataSheet.Range("H3").Value = 6

# Match the header/index-column styling used elsewhere in the workbook
# (same visual style -- bold, centered, thin-bordered -- as on "总计").
$totalSheet.Range("B1").Copy()
$q4DataSheet.Range("B1:H1").PasteSpecial(-4122)

$totalSheet.Range("A2").Copy()
$q4DataSheet.Range("A2:A3").PasteSpecial(-4122)
